$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: add columns D and E (copy formats from B and C) for rows 2-17 ---
$ws.Range("B2:B17").Copy() | Out-Null
$ws.Range("D2:D17").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:C17").Copy() | Out-Null
$ws.Range("E2:E17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Step 2: rewrite A/B/C values per row (re-sorted by region ascending) ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 10230.75
$ws.Range("C2").Value = 51153.75
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 33976.5
$ws.Range("C3").Value = 169882.5
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 16988.25
$ws.Range("C4").Value = 84941.25
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 23133
$ws.Range("C5").Value = 115665
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 47292
$ws.Range("C6").Value = 236460
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 36400.5
$ws.Range("C7").Value = 182002.5
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 51846
$ws.Range("C8").Value = 259230
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 49109.25
$ws.Range("C9").Value = 245546.25
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 31101
$ws.Range("C10").Value = 259230
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 29856.75
$ws.Range("C11").Value = 149283.75
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 4385.25
$ws.Range("C12").Value = 21926.25
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 6708.75
$ws.Range("C13").Value = 33543.75
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 141220.5
$ws.Range("C14").Value = 706102.5
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 15780.75
$ws.Range("C15").Value = 78903.75
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 7161.5249999999996
$ws.Range("C16").Value = 35807.625
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 19631.25
$ws.Range("C17").Value = 98156.25

# --- Step 3: add row 18 with SUM formula in B18 ---
$ws.Range("B18").Formula = "=SUM(B2:B17)"

# --- Step 4: sheet view / selection update ---
$ws.Range("E18").Select() | Out-Null
